$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Job adverts by occupation" latest period cell (note about job ads)
$ws.Range("D13").Value = "Sep 2025 (Oct 25)*"

# Scroll the view back so column A / row 1 is the top-left visible cell again,
# while keeping D14 as the selected (active) cell.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D14").Select()
